# Updates cryptos list values (Price column D, Volume(1h) column E)
# per commit "Updated cryptos list on Sun Sep 17 15:14:05 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are purely numeric-looking strings (e.g. "217.83").
# Excel auto-converts such text to a number on assignment, so force the
# affected cells to Text format first, then restore the default style so
# no unrelated formatting change is introduced.
$numericLookingCells = @("D5", "D11", "D16", "D19", "D21", "D22", "D24", "D25", "D29", "D30", "D33", "D35", "D38", "D39", "D41", "D43", "D45", "D46", "D48", "D50")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.750.01"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "1.640.96"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "217.83"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("E8").Value = "  +0.56%  "

$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("E10").Value = "  +0.42%  "

$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "1.869.33"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "1.657.33"
$ws.Range("E13").Value = "  +1.10%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D16").Value = "64.73"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "26.722.09"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  -0.93%  "

$ws.Range("D19").Value = "215.22"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("D22").Value = "2.39"
$ws.Range("E22").Value = "  +8.48%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "9.27"
$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").Value = "145.32"
$ws.Range("E25").Value = "  +0.48%  "

$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("D29").Value = "15.64"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  +0.76%  "

$ws.Range("D34").Value = "1.287.22"
$ws.Range("E34").Value = "  +1.17%  "

$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").Value = "0.538"
$ws.Range("E38").Value = "  +1.80%  "

$ws.Range("D39").Value = "0.818"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("E40").Value = "  +0.48%  "

$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("D43").Value = "5.28"
$ws.Range("E43").Value = "  -2.18%  "

$ws.Range("D44").Value = "1.779.95"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "61.17"
$ws.Range("E45").Value = "  +3.74%  "

$ws.Range("D46").Value = "91.84"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").Value = "0.0967"
$ws.Range("E50").Value = "  +0.69%  "

$ws.Range("E51").Value = "  -0.19%  "

# Restore default (Normal) style on the cells we temporarily reformatted,
# so only the cell content differs from the original workbook.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
